$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$sh = $s.Shapes.Item(16)

# Update the text from "Reorting Period" to "Measurement Period"
$sh.TextFrame.TextRange.Text = "Measurement Period"
$sh.TextFrame.TextRange.Runs(1).Font.Language = 1033

# Update position and size (EMU -> points, 1 pt = 12700 EMU)
$emuPerPt = 12700
$sh.Left = 4027223 / $emuPerPt
$sh.Top = 4255626 / $emuPerPt
$sh.Width = 2209064 / $emuPerPt
$sh.Height = 369332 / $emuPerPt
